$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B36").Value = "Data/HPA_tox_lists/bone+marrow_Tissue+enhanced.csv"
$ws.Range("B37").Value = "Data/HPA_tox_lists/bone+marrow_Tissue+enriched.csv"
$ws.Range("B38").Value = "Data/HPA_tox_lists/cervix_Group+enriched.csv"
$ws.Range("B39").Value = "Data/HPA_tox_lists/cervix_Tissue+enhanced.csv"
$ws.Range("B40").Value = "Data/HPA_tox_lists/cervix_Tissue+enriched.csv"
$ws.Range("B41").Value = "Data/HPA_tox_lists/heart_Group+enriched.csv"
$ws.Range("B42").Value = "Data/HPA_tox_lists/heart_Tissue+enhanced.csv"
$ws.Range("B43").Value = "Data/HPA_tox_lists/heart_Tissue+enriched.csv"
$ws.Range("B44").Value = "Data/HPA_tox_lists/kidney_Group+enriched.csv"
$ws.Range("B45").Value = "Data/HPA_tox_lists/kidney_Tissue+enhanced.csv"
$ws.Range("B46").Value = "Data/HPA_tox_lists/kidney_Tissue+enriched.csv"
$ws.Range("B47").Value = "Data/HPA_tox_lists/liver_Group+enriched.csv"
$ws.Range("B48").Value = "Data/HPA_tox_lists/liver_Tissue+enhanced.csv"
$ws.Range("B49").Value = "Data/HPA_tox_lists/liver_Tissue+enriched.csv"
$ws.Range("B50").Value = "Data/HPA_tox_lists/lymph+node_Group+enriched.csv"
$ws.Range("B51").Value = "Data/HPA_tox_lists/lymph+node_Tissue+enhanced.csv"
$ws.Range("B52").Value = "Data/HPA_tox_lists/lymph+node_Tissue+enriched.csv"
$ws.Range("B53").Value = "Data/HPA_tox_lists/ovary_Group+enriched.csv"
$ws.Range("B54").Value = "Data/HPA_tox_lists/ovary_Tissue+enhanced.csv"
$ws.Range("B55").Value = "Data/HPA_tox_lists/ovary_Tissue+enriched.csv"
$ws.Range("B56").Value = "Data/HPA_tox_lists/testis_Group+enriched.csv"
$ws.Range("B57").Value = "Data/HPA_tox_lists/testis_Tissue+enhanced.csv"
$ws.Range("B58").Value = "Data/HPA_tox_lists/testis_Tissue+enriched.csv"
$ws.Range("B35").Value = "Data/HPA_tox_lists/bone+marrow_Group+enriched.csv"

$ws.Range("A35").Value = "Human Protein Atlas: RNA expression in Bone marrow at Group enriched level"
$ws.Range("A36").Value = "Human Protein Atlas: RNA expression in Bone marrow at Tissue enhanced level"
$ws.Range("A37").Value = "Human Protein Atlas: RNA expression in Bone marrow at Tissue enriched level"
$ws.Range("A38").Value = "Human Protein Atlas: RNA expression in Cervix at Group enriched level"
$ws.Range("A39").Value = "Human Protein Atlas: RNA expression in Cervix at Tissue enhanced level"
$ws.Range("A40").Value = "Human Protein Atlas: RNA expression in Cervix at Tissue enriched level"
$ws.Range("A41").Value = "Human Protein Atlas: RNA expression in Heart at Group enriched level"
$ws.Range("A42").Value = "Human Protein Atlas: RNA expression in Heart at Tissue enhanced level"
$ws.Range("A43").Value = "Human Protein Atlas: RNA expression in Heart at Tissue enriched level"
$ws.Range("A44").Value = "Human Protein Atlas: RNA expression in Kidney at Group enriched level"
$ws.Range("A45").Value = "Human Protein Atlas: RNA expression in Kidney at Tissue enhanced level"
$ws.Range("A46").Value = "Human Protein Atlas: RNA expression in Kidney at Tissue enriched level"
$ws.Range("A47").Value = "Human Protein Atlas: RNA expression in Liver at Group enriched level"
$ws.Range("A48").Value = "Human Protein Atlas: RNA expression in Liver at Tissue enhanced level"
$ws.Range("A49").Value = "Human Protein Atlas: RNA expression in Liver at Tissue enriched level"
$ws.Range("A50").Value = "Human Protein Atlas: RNA expression in Lymph node at Group enriched level"
$ws.Range("A51").Value = "Human Protein Atlas: RNA expression in Lymph node at Tissue enhanced level"
$ws.Range("A52").Value = "Human Protein Atlas: RNA expression in Lymph node at Tissue enriched level"
$ws.Range("A53").Value = "Human Protein Atlas: RNA expression in Ovary at Group enriched level"
$ws.Range("A54").Value = "Human Protein Atlas: RNA expression in Ovary at Tissue enhanced level"
$ws.Range("A55").Value = "Human Protein Atlas: RNA expression in Ovary at Tissue enriched level"
$ws.Range("A56").Value = "Human Protein Atlas: RNA expression in Testis at Group enriched level"
$ws.Range("A57").Value = "Human Protein Atlas: RNA expression in Testis at Tissue enhanced level"
$ws.Range("A58").Value = "Human Protein Atlas: RNA expression in Testis at Tissue enriched level"

$ws.Range("A50").Select()
